$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 16 and Row 18 swap their worker data (doc number, name, period, valor mora);
# Row 17 (SANDRA MARIA PANIZA ARIZA / 22785827 / 2004) is unchanged.

# Row 16: now MANUEL ANTONIO ARIZA SALGADO / 1051888181 / 2001 / 28090
$ws.Range("C16").Value = "1051888181"
$ws.Range("D16").Value = "MANUEL ANTONIO ARIZA SALGADO"
$ws.Range("E16").Value = "2001"
$ws.Range("F16").Value = 28090

# Row 18: now SANDRA MARIA PANIZA ARIZA / 22785827 / 2005 / 35112
$ws.Range("C18").Value = "22785827"
$ws.Range("D18").Value = "SANDRA MARIA PANIZA ARIZA"
$ws.Range("E18").Value = "2005"
$ws.Range("F18").Value = 35112
